# Update "想去人数" (want-to-go count) figures in column F for the
# exhibition sheets, matching the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 603
    6  = 182
    7  = 63
    10 = 5548
    11 = 4928
    13 = 42
    16 = 203
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
